$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-43 is being bumped by one day
# (Excel serial date 45756 -> 45757, i.e. 2025-04-09 -> 2025-04-10).
for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45756) {
        $cell.Value2 = 45757
    }
}
